$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "UEGC" ticker is being removed from the Stocks list (column B). The
# numeric index in column A is a plain 0-based row counter and is left as-is;
# only the tickers in column B shift up to fill the gap, and the now-unused
# last row is removed entirely.
$find = $ws.Range("B:B").Find("UEGC")
$startRow = $find.Row

# Find the last used row in column B.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Shift tickers in column B up by one, starting at the row that held "UEGC".
for ($r = $startRow; $r -lt $lastRow; $r++) {
    $src = $ws.Cells.Item($r + 1, 2)
    $dst = $ws.Cells.Item($r, 2)
    $dst.Value = $src.Value2
}

# Remove the now-duplicated last row entirely (shrinks the used range).
$ws.Rows($lastRow).Delete()
